# petty-cashBook-2021.xlsx -- "Update 7-Jul-2021, end of day update."
# Edits the "Sheet1" worksheet (tab-selected, the daily "Buku KAS HARIAN"-style
# petty cash ledger that runs A1:L114) to record the 6/Jul and 7/Jul entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11: Wages Expense for 6/Jul grew by an extra 260,000 ---
$ws.Range("D11").Formula = "=60000+260000"

# --- Row 14: freight/claim total grew (632000 -> 642000, plus a new 4,429,000) ---
$ws.Range("D14").Formula = "=4500000+432000+642000+4429000"

# --- New entries, rows 16-24 (previously-blank trailer rows) ---

# Row 16: A/R
$ws.Range("B16").Value = "A/R"
$ws.Range("C16").Formula = "=12960000"

# Row 17: SALES - cash/retail
$ws.Range("B17").Value = "SALES - cash/retail"
$ws.Range("C17").Formula = "=6498725+10928275-12960000"

# Row 18: SELISIH - lebih
$ws.Range("B18").Value = "SELISIH - lebih"
$ws.Range("C18").Value = 460000

# Row 19: SETOR KE BANK
$ws.Range("B19").Value = "SETOR KE BANK"
$ws.Range("D19").Value = 6000000

# Row 20: 7/Jul/2021, Wages Expense
$ws.Range("A20").Value = 44384
$ws.Range("B20").Value = "Wages Expense"
$ws.Range("D20").Formula = "=60000"

# Row 21: TRANSFER BCA
$ws.Range("B21").Value = "TRANSFER BCA"
$ws.Range("D21").Formula = "=1405000+1864000+6027200+9027000+3000000+9027000+2606000+1897000"

# Row 22: STNK - suzuki (new description, not previously in the shared strings)
$ws.Range("B22").Value = "STNK - suzuki"
$ws.Range("D22").Formula = "=529000"

# Row 23: A/R
$ws.Range("B23").Value = "A/R"
$ws.Range("C23").Formula = "=6027200+9027000+3000000+9027000+29760000+2281500+1897000"

# Row 24: FREIGHT - OUT
$ws.Range("B24").Value = "FREIGHT - OUT"
$ws.Range("D24").Formula = "=60000"

# --- View state: move the cursor down to where today's entries land ---
$ws.Range("D40").Select()
